$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "test"
$ws.Range("B6").Value = 1934
$ws.Range("C6").Value = 1661
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 1633
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "36.96"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "31.75"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "0.08"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "31.21"
$ws.Range("I6").ClearFormats()
$ws.Range("J6").Value = "2025-08-20 18:06:10"

$ws.Range("A7").Value = "test"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 3599
$ws.Range("E7").Value = 1633
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "0.00"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.00"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "68.79"
$ws.Range("H7").ClearFormats()
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "31.21"
$ws.Range("I7").ClearFormats()
$ws.Range("J7").Value = "2025-08-20 22:30:11"

$ws.Range("A8").Value = "test"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 3599
$ws.Range("E8").Value = 1633
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0.00"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.00"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "68.79"
$ws.Range("H8").ClearFormats()
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "31.21"
$ws.Range("I8").ClearFormats()
$ws.Range("J8").Value = "2025-08-20 22:32:14"

$ws.Range("A9").Value = "test"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 3599
$ws.Range("E9").Value = 1633
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "0.00"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.00"
$ws.Range("G9").ClearFormats()
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "68.79"
$ws.Range("H9").ClearFormats()
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "31.21"
$ws.Range("I9").ClearFormats()
$ws.Range("J9").Value = "2025-08-20 22:33:59"

$ws.Range("A10").Value = "test"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 3599
$ws.Range("E10").Value = 1633
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "0.00"
$ws.Range("F10").ClearFormats()
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.00"
$ws.Range("G10").ClearFormats()
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "68.79"
$ws.Range("H10").ClearFormats()
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "31.21"
$ws.Range("I10").ClearFormats()
$ws.Range("J10").Value = "2025-08-20 22:34:45"

$ws.Range("A11").Value = "test"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 3599
$ws.Range("E11").Value = 1633
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "0.00"
$ws.Range("F11").ClearFormats()
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.00"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "68.79"
$ws.Range("H11").ClearFormats()
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "31.21"
$ws.Range("I11").ClearFormats()
$ws.Range("J11").Value = "2025-08-20 22:35:06"

$ws.Range("A12").Value = "test"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 3599
$ws.Range("E12").Value = 1633
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "0.00"
$ws.Range("F12").ClearFormats()
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.00"
$ws.Range("G12").ClearFormats()
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "68.79"
$ws.Range("H12").ClearFormats()
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "31.21"
$ws.Range("I12").ClearFormats()
$ws.Range("J12").Value = "2025-08-20 22:39:48"

$ws.Range("A13").Value = "test"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 3599
$ws.Range("E13").Value = 1633
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "0.00"
$ws.Range("F13").ClearFormats()
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.00"
$ws.Range("G13").ClearFormats()
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "68.79"
$ws.Range("H13").ClearFormats()
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "31.21"
$ws.Range("I13").ClearFormats()
$ws.Range("J13").Value = "2025-08-20 22:41:41"

$ws.Range("A14").Value = "test"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 10789
$ws.Range("E14").Value = 1633
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0.00"
$ws.Range("F14").ClearFormats()
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.00"
$ws.Range("G14").ClearFormats()
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "86.85"
$ws.Range("H14").ClearFormats()
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "13.15"
$ws.Range("I14").ClearFormats()
$ws.Range("J14").Value = "2025-08-20 22:42:12"

$ws.Range("A15").Value = "test"
$ws.Range("B15").Value = 1934
$ws.Range("C15").Value = 1661
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 1633
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "36.96"
$ws.Range("F15").ClearFormats()
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "31.75"
$ws.Range("G15").ClearFormats()
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "0.08"
$ws.Range("H15").ClearFormats()
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "31.21"
$ws.Range("I15").ClearFormats()
$ws.Range("J15").Value = "2025-08-20 23:12:30"

$ws.Range("A16").Value = "test"
$ws.Range("B16").Value = 1934
$ws.Range("C16").Value = 1661
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 1633
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "36.96"
$ws.Range("F16").ClearFormats()
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "31.75"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "0.08"
$ws.Range("H16").ClearFormats()
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "31.21"
$ws.Range("I16").ClearFormats()
$ws.Range("J16").Value = "2025-08-20 23:13:13"

